$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, $tr1.Length).Text = "Slide 1"

$tr2 = $s1.Shapes.Item(3).TextFrame.TextRange
$tr2.Characters(1, $tr2.Length).Text = "an image"

$s2 = $p.Slides.Item(2)
$tr3 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, $tr3.Length).Text = "Slide 2"

$tr4 = $s2.Shapes.Item(4).TextFrame.TextRange
$tr4.Characters(1, $tr4.Length).Text = "an image"
